$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("F41:F42").Formula = "=1-D41/E41"

$ws = $wb.Worksheets.Item("Intercooler number")
$ws.Rows.Item(15).Insert()
$ws.Range("C15").Formula = '=1-C14/$E$14'
$ws.Range("D15").Formula = '=1-D14/$E$14'
$ws.Range("E15").Formula = '=1-E14/$E$14'
$ws.Range("F15:I15").Formula = '=1-F14/$E$14'
$ws.Range("C21:I21").Formula = '=$B$20-C20'
$ws.Range("B3").Select() | Out-Null
$ws.Activate() | Out-Null
